# Append summary rows/formulas to Sheet1, matching the target commit
# ("xlxs for freelancer, k=0.1, k=2"): an AVERAGE check in J12, plus four
# labeled summary statistics (average/worst of the SW and SC ratios) in
# A14:B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: sanity-check average of the k column (J) -----------------
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# --- Rows 14-17: summary labels (A) + stat formulas (B) ----------------
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$statRange = $ws.Range("B14:B17")
$statRange.Font.Bold = $true
$statRange.Font.Size = 12
$statRange.VerticalAlignment = -4108

# --- Page setup / view cosmetics, mirroring the re-saved workbook ------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("J12").Select()
